# Update the cheat sheet (man/figures/cheatsheet/cheatsheet.pptx)
#
# Applies four text edits on the single slide of the deck:
#   1. Footer shape: package version "0.5.3" -> "0.5.4"
#   2. Footer shape: updated date "2021-05" -> "2021-06"
#   3. Body shape: "...based on the Poisson, binomial, or hypergeometric
#      likelihood..." -> "...based on the binomial, Poisson, or
#      hypergeometric likelihood..."
#   4. Body shape: "...statistical results their interpretation." ->
#      "...statistical results and their interpretation."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# --- Shape id 322: footer line with package version / update date ---
$shFooter = Get-ShapeById $s 322
$trFooter = $shFooter.TextFrame.TextRange
$fullFooter = $trFooter.Text

$oldVersion = "0.5.3"
$newVersionDigit = "4"
$idxVersion = $fullFooter.IndexOf($oldVersion)
if ($idxVersion -ge 0) {
    # last character of "0.5.3" is the "3" -> replace with "4"
    $cVersion = $trFooter.Characters($idxVersion + $oldVersion.Length, 1)
    $cVersion.Text = $newVersionDigit
}

# Re-read text since it changed above
$fullFooter = $trFooter.Text
$oldMonth = "2021-05"
$newMonth = "06"
$idxMonth = $fullFooter.IndexOf($oldMonth)
if ($idxMonth -ge 0) {
    # last two characters of "2021-05" are the "05" -> replace with "06"
    $cMonth = $trFooter.Characters($idxMonth + $oldMonth.Length - 1, 2)
    $cMonth.Text = $newMonth
}

# --- Shape id 124: distribution wording ---
$shDist = Get-ShapeById $s 124
$trDist = $shDist.TextFrame.TextRange
$fullDist = $trDist.Text

$anchor = "Poisson, binomial, or hypergeometric"
$idxAnchor = $fullDist.IndexOf($anchor)
if ($idxAnchor -ge 0) {
    # "P" is its own run; replace it with "binomial"
    $cP = $trDist.Characters($idxAnchor + 1, 1)
    # "oisson, binomial, or hypergeometric " (36 chars incl. trailing space) is the next run
    $cRest = $trDist.Characters($idxAnchor + 2, 36)
    $cRest.Text = ", Poisson, or hypergeometric "
    $cP.Text = "binomial"
}

# --- Shape id 139: add missing "and" ---
$shReport = Get-ShapeById $s 139
$trReport = $shReport.TextFrame.TextRange
$fullReport = $trReport.Text

$oldTail = " and creates a report containing the statistical results their interpretation."
$newTail = " and creates a report containing the statistical results and their interpretation."
$idxTail = $fullReport.IndexOf($oldTail)
if ($idxTail -ge 0) {
    $cTail = $trReport.Characters($idxTail + 1, $oldTail.Length)
    $cTail.Text = $newTail
}

Write-Host "Footer text now:" $trFooter.Text
Write-Host "Distribution text now:" $trDist.Text
Write-Host "Report text now:" $trReport.Text
